$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the start time / end time values (stored as raw nanosecond timestamps).
$ws.Range("B1").Value = [double]"1.5108736928179999E+18"
$ws.Range("B2").Value = [double]"1.510873839248E+18"

# Formulas in B3 (=B2-B1) and B4 (=B3/POWER(10,9)) recalculate automatically.
$excel.Calculate()
